$d = $word.ActiveDocument

function Set-RunFont($searchText) {
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $f.Text = $searchText
    $f.Replacement.Text = $searchText
    $f.Replacement.Font.Name = "Times New Roman"
    $f.Replacement.Font.NameFarEast = "Times New Roman"
    $f.Replacement.Font.NameBi = "Times New Roman"
    $f.Replacement.Font.NameOther = "Times New Roman"
    $f.Replacement.Font.Size = 12
    $result = $f.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
    return $result
}

# --- 1) Committee info table: give every run in the cell explicit
#        Times New Roman / 12pt formatting (keeping existing italics). ---
$committeeRuns = @(
    "2.1 Председатель государственной экзаменационной комиссии: ",
    "Алексеев Алексей Алексеевич, ",
    "информация о Алексееве Алексее Алексеевиче,",
    "утвержден приказом от 17.12.2025 1234/1 (с изменениями и дополнениями).",
    "Электронный адрес: alekseev@gmail.com",
    "Контактный телефон: +72222222222",
    "2.2 ",
    "Романов Роман Романович, ",
    "информация о Романове Романе Романовиче.",
    "Электронный адрес: romanov@gmail.com",
    "Контактный телефон: +71111111111"
)

foreach ($t in $committeeRuns) {
    Set-RunFont $t | Out-Null
}

# --- 2) Collapse the "Время начала заседания 14:30." run-chain into a
#        single run that reads "Время начала заседания: _____ ". ---
$rng2 = $d.Content
$f2 = $rng2.Find
$f2.ClearFormatting()
$f2.Replacement.ClearFormatting()
$searchText2 = "Время начала заседания 14:30."
$replaceText2 = "Время начала заседания: _____ "
$f2.Text = $searchText2
$f2.Replacement.Text = $replaceText2
$f2.Execute($searchText2, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText2, 2) | Out-Null
